$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-enter the "1 - One" value in G2 with a leading apostrophe (quote prefix),
# matching the text-forced formatting applied by Excel.
$ws.Range("G2").Value = "'1 - One"

# Update the check-in / check-out dates, forcing them to be stored as text
# (number format "@") since the new values can't be parsed as valid dates.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "28/02/2023"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "27/02/2023"

# Update the active selection on the sheet to I13.
[void]$ws.Range("I13").Select()
